# Add a new "allowsplits" parameter row to the "parameter" worksheet,
# right after the existing "plan.webservice" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parameter")

$ws.Cells.Item(11, 1).Value = "allowsplits"
$ws.Cells.Item(11, 2).Value = "'true"
$ws.Cells.Item(11, 2).Style = "Normal"
$ws.Cells.Item(11, 3).Value = "Controls whether a sales order or forecast can be split across multiple manufacturing orders during planning. Default: false"

$ws.Select()
$ws.Range("C11").Select()
